$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $val)
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "42.997.78"
Set-TextValue "E2" "  +1.93%  "
# Row 3
Set-TextValue "D3" "2.310.91"
Set-TextValue "E3" "  +1.90%  "
# Row 4
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.00%  "
# Row 5
Set-TextValue "D5" "303.56"
Set-TextValue "E5" "  +1.52%  "
# Row 6
Set-TextValue "D6" "101.15"
Set-TextValue "E6" "  +5.87%  "
# Row 7
Set-TextValue "E7" "  +1.81%  "
# Row 9
Set-TextValue "E9" "  +4.21%  "
# Row 10
Set-TextValue "D10" "35.13"
Set-TextValue "E10" "  +5.32%  "
# Row 11
Set-TextValue "D11" "0.0796"
Set-TextValue "E11" "  +0.81%  "
# Row 12
Set-TextValue "E12" "  +4.07%  "
# Row 13
Set-TextValue "D13" "17.89"
Set-TextValue "E13" "  +15.21%  "
# Row 14
Set-TextValue "D14" "6.92"
Set-TextValue "E14" "  +4.06%  "
# Row 15
Set-TextValue "D15" "2.686.31"
Set-TextValue "E15" "  +2.38%  "
# Row 16
Set-TextValue "D16" "2.282.17"
Set-TextValue "E16" "  -0.02%  "
# Row 17
Set-TextValue "E17" "  +4.22%  "
# Row 18
Set-TextValue "D18" "42.946.67"
Set-TextValue "E18" "  +1.97%  "
# Row 19
Set-TextValue "D19" "12.65"
Set-TextValue "E19" "  +8.07%  "
# Row 20
Set-TextValue "D20" "6.17"
Set-TextValue "E20" "  +3.19%  "
# Row 21
Set-TextValue "D21" "0.0₃0904"
Set-TextValue "E21" "  +1.52%  "
# Row 22
Set-TextValue "E22" "  +1.92%  "
# Row 23
Set-TextValue "D23" "237.70"
Set-TextValue "E23" "  +1.27%  "
# Row 24
Set-TextValue "E24" "  +12.95%  "
# Row 25
Set-TextValue "D25" "2.48"
Set-TextValue "E25" "  +1.00%  "
# Row 26
Set-TextValue "E26" "  -0.12%  "
# Row 27
Set-TextValue "D27" "24.74"
Set-TextValue "E27" "  +3.14%  "
# Row 28
Set-TextValue "D28" "2.24"
Set-TextValue "E28" "  -2.02%  "
# Row 29
Set-TextValue "D29" "167.67"
Set-TextValue "E29" "  -0.38%  "
# Row 30
Set-TextValue "B30" "Cosmos"
Set-TextValue "C30" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D30" "9.23"
Set-TextValue "E30" "  +0.92%  "
# Row 31
Set-TextValue "B31" "InjectiveProtocol"
Set-TextValue "C31" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D31" "33.99"
Set-TextValue "E31" "  +0.36%  "
# Row 32
Set-TextValue "E32" "  +0.05%  "
# Row 33
Set-TextValue "D33" "5.01"
Set-TextValue "E33" "  +2.22%  "
# Row 34
Set-TextValue "D34" "4.62"
Set-TextValue "E34" "  +3.16%  "
# Row 35
Set-TextValue "E35" "  +3.59%  "
# Row 36
Set-TextValue "D36" "17.07"
Set-TextValue "E36" "  +3.21%  "
# Row 37
Set-TextValue "D37" "0.0693"
Set-TextValue "E37" "  +1.00%  "
# Row 38
Set-TextValue "E38" "  +3.78%  "
# Row 39
Set-TextValue "D39" "1.79"
Set-TextValue "E39" "  +4.02%  "
# Row 40
Set-TextValue "E40" "  +1.63%  "
# Row 41
Set-TextValue "E41" "  +0.79%  "
# Row 42
Set-TextValue "D42" "2.002.84"
Set-TextValue "E42" "  +2.24%  "
# Row 43
Set-TextValue "D43" "2.27"
Set-TextValue "E43" "  -6.46%  "
# Row 44
Set-TextValue "E44" "  +3.77%  "
# Row 45
Set-TextValue "D45" "10.24"
Set-TextValue "E45" "  +7.41%  "
# Row 46
Set-TextValue "D46" "17.57"
Set-TextValue "E46" "  +0.97%  "
# Row 47
Set-TextValue "D47" "2.84"
Set-TextValue "E47" "  +2.28%  "
# Row 48
Set-TextValue "D48" "55.43"
Set-TextValue "E48" "  +5.91%  "
# Row 49
Set-TextValue "D49" "2.529.50"
Set-TextValue "E49" "  +1.39%  "
# Row 50
Set-TextValue "D50" "1.53"
Set-TextValue "E50" "  +4.55%  "
# Row 51
Set-TextValue "D51" "4.56"
Set-TextValue "E51" "  +0.59%  "
